# Add 2022-Q4 data:
#  - Insert a new "2022-Q4" detail worksheet between "总计" and "2022-Q3".
#  - Insert a new summary row for 2022-Q4 at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and
#    push the existing 2022-Q3 summary row down to row 3.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Copy the (now shifted) row 3 formatting into the new row 2 so both rows
# share the same look (bold index column, etc.)
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 26
$totalSheet.Range("D2").Value = 10.84

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 16
$totalSheet.Range("D3").Value = 7.84

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q4" detail worksheet, positioned right after
#    "总计" and before "2022-Q3".
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($j = 0; $j -lt $headers.Count; $j++) {
    $q4Sheet.Cells.Item(1, $j + 2).Value = $headers[$j]
}

$q4data = @(
  @(0, "001071", "华安媒体互联网混合A", "51.25", "91.12", "4.08", "2.0910", 2),
  @(1, "004666", "长城久嘉创新成长灵活配置混合A", "24.38", "92.75", "8.10", "1.9748", 1),
  @(2, "001694", "华安沪港深外延增长混合A", "38.21", "94.15", "3.69", "1.4099", 5),
  @(3, "010052", "长城久嘉创新成长灵活配置混合C", "15.65", "92.75", "8.10", "1.2676", 1),
  @(4, "006879", "华安智能生活混合A", "28.87", "92.19", "3.53", "1.0191", 5),
  @(5, "007460", "华安成长创新混合A", "16.41", "93.15", "4.08", "0.6695", 5),
  @(6, "013621", "华安智能生活混合C", "11.19", "92.19", "3.53", "0.3950", 5),
  @(7, "002621", "中欧消费主题股票A", "10.21", "94.00", "3.68", "0.3757", 5),
  @(8, "008980", "中邮科技创新精选混合A", "8.35", "88.21", "3.39", "0.2831", 8),
  @(9, "007126", "博道远航混合A", "6.63", "94.00", "3.68", "0.2440", 5),
  @(10, "002697", "中欧消费主题股票C", "5.42", "94.00", "3.68", "0.1995", 5),
  @(11, "014754", "华安景气优选混合A", "5.19", "92.65", "3.64", "0.1889", 5),
  @(12, "007127", "博道远航混合C", "4.92", "94.00", "3.68", "0.1811", 5),
  @(13, "013620", "华安媒体互联网混合C", "2.39", "91.12", "4.08", "0.0975", 2),
  @(14, "008981", "中邮科技创新精选混合C", "2.29", "88.21", "3.39", "0.0776", 8),
  @(15, "006025", "诺安优化配置混合", "0.77", "91.43", "9.64", "0.0742", 2),
  @(16, "013369", "汇添富自主核心科技一年持有混合A", "2.40", "80.35", "2.74", "0.0658", 8),
  @(17, "014755", "华安景气优选混合C", "1.48", "92.65", "3.64", "0.0539", 5),
  @(18, "013340", "创金合信芯片产业股票C", "0.94", "92.41", "4.84", "0.0455", 6),
  @(19, "013339", "创金合信芯片产业股票A", "0.92", "92.41", "4.84", "0.0445", 6),
  @(20, "001662", "创金沪港深精选混合", "0.67", "93.27", "2.90", "0.0194", 8),
  @(21, "016099", "华安成长创新混合C", "0.47", "93.15", "4.08", "0.0192", 5),
  @(22, "013370", "汇添富自主核心科技一年持有混合C", "0.66", "80.35", "2.74", "0.0181", 8),
  @(23, "015919", "申万菱信专精特新主题混合A", "0.39", "48.55", "4.46", "0.0174", 3),
  @(24, "014972", "华安沪港深外延增长混合C", "0.13", "94.15", "3.69", "0.0048", 5),
  @(25, "015920", "申万菱信专精特新主题混合C", "0.02", "48.55", "4.46", "0.0009", 3)
)

for ($i = 0; $i -lt $q4data.Count; $i++) {
    $r = $i + 2
    $row = $q4data[$i]

    # Force the fund-code / numeric-looking text columns (B..G) to remain
    # plain text so things like leading zeros ("001071") and percentages
    # are preserved instead of being auto-converted to numbers.
    $q4Sheet.Range("B" + $r + ":G" + $r).NumberFormat = "@"

    $q4Sheet.Cells.Item($r, 1).Value = $row[0]
    $q4Sheet.Cells.Item($r, 2).Value = $row[1]
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]
    $q4Sheet.Cells.Item($r, 4).Value = $row[3]
    $q4Sheet.Cells.Item($r, 5).Value = $row[4]
    $q4Sheet.Cells.Item($r, 6).Value = $row[5]
    $q4Sheet.Cells.Item($r, 7).Value = $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]
}

# Match formatting with the rest of the workbook: bold header row and bold
# index column, same style as used on the "总计" sheet.
$totalSheet.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Restore the original active sheet (the "总计" sheet was the active tab).
# ---------------------------------------------------------------------------
$totalSheet.Activate()
$totalSheet.Range("A1").Select()
